$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-205) holds the "Förändrad" date, stored as serial date 45189.
# Update it to 45190 for every data row.
$range = $ws.Range("C2:C205")
$range.Value = 45190
